$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A6").Value = "07 Jan 2019"
$ws.Range("B6").Value = 0
